# Article 96 is live
# Shift the rotating "blog" slots in row 7 down by one, dropping the
# oldest (ser: 93) and introducing the newly published article (ser: 96)
# into the slot vacated at the front of the rotation (C7).
#
#   I7: ser: 93 -> ser: 94   (was E7's old value)
#   E7: ser: 94 -> ser: 95   (was C7's old value)
#   C7: ser: 95 -> ser: 96   (brand new article text)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newArticle = "type: blog" + [char]10 + "width: 2" + [char]10 + "height: 1" + [char]10 + "ser: 96"

$ws.Range("I7").Value = $ws.Range("E7").Value()
$ws.Range("E7").Value = $ws.Range("C7").Value()
$ws.Range("C7").Value = $newArticle

$ws.Range("I7").Select()
